$wb = $excel.ActiveWorkbook


# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1111398.1
$ws.Range("I6").Value = 1250285.4
$ws.Range("K6").Value = 3750856.2
$ws.Range("M6").Value = -3750744.2
$ws.Range("H8").Value = 48.666668
$ws.Range("I8").Value = 48.666668
$ws.Range("K8").Value = 146.000004
$ws.Range("M8").Value = -7.00000399999999
$ws.Range("H11").Value = 344.69232
$ws.Range("I11").Value = 344.69232
$ws.Range("K11").Value = 344.69232
$ws.Range("M11").Value = -204.69232
$ws.Range("H31").Value = 74
$ws.Range("I31").Value = 74
$ws.Range("K31").Value = 222
$ws.Range("M31").Value = 8
$ws.Range("H33").Value = 7143659
$ws.Range("I33").Value = 12500655
$ws.Range("K33").Value = 12500655
$ws.Range("M33").Value = -12500426
$ws.Range("H40").Value = 2399
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H41").Value = 2581
$ws.Range("I41").Value = 3209.3635
$ws.Range("K41").Value = 3209.3635
$ws.Range("M41").Value = -2769.3635
$ws.Range("H51").Value = 3117.647
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 2500
$ws.Range("M51").Value = -2016
$ws.Range("H58").Value = 181.09091
$ws.Range("I58").Value = 181.09091
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 543.27273
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -393.27273
$ws.Range("N58").ClearContents()
$ws.Range("H92").Value = 2158.6
$ws.Range("I92").Value = 1597.6666
$ws.Range("K92").Value = 1597.6666
$ws.Range("M92").Value = -349.6666
$ws.Range("H98").Value = 3187.5
$ws.Range("I98").Value = 1375
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 1375
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = 123
$ws.Range("N98").Value = -7996
$ws.Range("H106").Value = 3726.6428
$ws.Range("I106").Value = 3709.75
$ws.Range("K106").Value = 3709.75
$ws.Range("M106").Value = -3078.75
$ws.Range("H122").Value = 3187.5
$ws.Range("I122").Value = 1375
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4125
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1675
$ws.Range("N122").Value = -19900
$ws.Range("H125").Value = 2437.5
$ws.Range("I125").Value = 2437.5
$ws.Range("K125").Value = 21937.5
$ws.Range("M125").Value = -19477.5
$ws.Range("H133").Value = 70661.11
$ws.Range("J133").Value = 54493.75
$ws.Range("L133").Value = 54493.75
$ws.Range("N133").Value = -64613.75
$ws.Range("H138").Value = 3602.2812
$ws.Range("I138").Value = 1627.3334
$ws.Range("J138").Value = 4155.2666
$ws.Range("K138").Value = 4882.0002
$ws.Range("L138").Value = 12465.7998
$ws.Range("M138").Value = 257.9997999999996
$ws.Range("N138").Value = -22745.7998
$ws.Range("H141").Value = 3252.9656
$ws.Range("I141").Value = 3341.84
$ws.Range("J141").Value = 2697.5
$ws.Range("K141").Value = 10025.52
$ws.Range("L141").Value = 8092.5
$ws.Range("M141").Value = -4845.52
$ws.Range("N141").Value = -18452.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 252
$ws.Range("I4").Value = 227.5
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 227.5
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = -111.5
$ws.Range("N4").Value = -582
$ws.Range("H28").Value = 2819.6
$ws.Range("I28").Value = 2819.6
$ws.Range("K28").Value = 2819.6
$ws.Range("M28").Value = -2627.6
$ws.Range("H31").Value = 4995.1665
$ws.Range("I31").Value = 4995.1665
$ws.Range("K31").Value = 4995.1665
$ws.Range("M31").Value = -4701.1665
$ws.Range("H61").Value = 4133.375
$ws.Range("I61").Value = 3942.2666
$ws.Range("K61").Value = 3942.2666
$ws.Range("M61").Value = -3730.2666
$ws.Range("H63").Value = 2963.2856
$ws.Range("I63").Value = 2949.0557
$ws.Range("K63").Value = 2949.0557
$ws.Range("M63").Value = -2263.0557
$ws.Range("H66").Value = 2963.2856
$ws.Range("I66").Value = 2949.0557
$ws.Range("K66").Value = 14745.2785
$ws.Range("M66").Value = -11313.2785
$ws.Range("H74").Value = 21744800
$ws.Range("I74").Value = 1092.625
$ws.Range("K74").Value = 1092.625
$ws.Range("M74").Value = -218.625
$ws.Range("H77").Value = 21744800
$ws.Range("I77").Value = 1092.625
$ws.Range("K77").Value = 5463.125
$ws.Range("M77").Value = -1095.125
$ws.Range("H88").Value = 2310.9565
$ws.Range("I88").Value = 2535.1538
$ws.Range("J88").Value = 2019.5
$ws.Range("K88").Value = 2535.1538
$ws.Range("L88").Value = 2019.5
$ws.Range("M88").Value = -2129.1538
$ws.Range("N88").Value = -2831.5
$ws.Range("H91").Value = 2310.9565
$ws.Range("I91").Value = 2535.1538
$ws.Range("J91").Value = 2019.5
$ws.Range("K91").Value = 2535.1538
$ws.Range("L91").Value = 2019.5
$ws.Range("M91").Value = -1131.1538
$ws.Range("N91").Value = -4827.5
$ws.Range("H99").Value = 2819.6
$ws.Range("I99").Value = 2819.6
$ws.Range("K99").Value = 2819.6
$ws.Range("M99").Value = 175.4000000000001
$ws.Range("H102").Value = 6014.5713
$ws.Range("I102").Value = 6014.5713
$ws.Range("K102").Value = 6014.5713
$ws.Range("M102").Value = -4392.5713
$ws.Range("H109").Value = 37474.5
$ws.Range("J109").Value = 34949
$ws.Range("L109").Value = 34949
$ws.Range("N109").Value = -37723
$ws.Range("H122").Value = 3407.1304
$ws.Range("I122").Value = 3256.2354
$ws.Range("J122").Value = 3834.6667
$ws.Range("K122").Value = 9768.706200000001
$ws.Range("L122").Value = 11504.0001
$ws.Range("M122").Value = -7318.706200000001
$ws.Range("N122").Value = -16404.0001
$ws.Range("H132").Value = 3012.8215
$ws.Range("I132").Value = 2642.36
$ws.Range("J132").Value = 6100
$ws.Range("K132").Value = 7927.08
$ws.Range("L132").Value = 18300
$ws.Range("M132").Value = -5397.08
$ws.Range("N132").Value = -23360
$ws.Range("H136").Value = 4133.375
$ws.Range("I136").Value = 3942.2666
$ws.Range("K136").Value = 11826.7998
$ws.Range("M136").Value = -9276.799800000001
$ws.Range("H139").Value = 75492.5
$ws.Range("J139").Value = 75492.5
$ws.Range("L139").Value = 75492.5
$ws.Range("N139").Value = -85772.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 81624.5
$ws.Range("J57").Value = 81624.5
$ws.Range("L57").Value = 81624.5
$ws.Range("N57").Value = -83064.5
$ws.Range("H86").Value = 4118.8076
$ws.Range("I86").Value = 3016.0667
$ws.Range("J86").Value = 5622.5454
$ws.Range("K86").Value = 3016.0667
$ws.Range("L86").Value = 5622.5454
$ws.Range("M86").Value = -1893.0667
$ws.Range("N86").Value = -7868.5454
$ws.Range("H89").Value = 4118.8076
$ws.Range("I89").Value = 3016.0667
$ws.Range("J89").Value = 5622.5454
$ws.Range("K89").Value = 15080.3335
$ws.Range("L89").Value = 28112.727
$ws.Range("M89").Value = -9464.333499999999
$ws.Range("N89").Value = -39344.727
$ws.Range("H99").Value = 5386.2
$ws.Range("I99").Value = 5386.2
$ws.Range("K99").Value = 5386.2
$ws.Range("M99").Value = -3888.2
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H134").Value = 1595.7297
$ws.Range("I134").Value = 1189.1818
$ws.Range("J134").Value = 4949.75
$ws.Range("K134").Value = 3567.5454
$ws.Range("L134").Value = 14849.25
$ws.Range("M134").Value = -1032.5454
$ws.Range("N134").Value = -19919.25
$ws.Range("H136").Value = 81624.5
$ws.Range("J136").Value = 81624.5
$ws.Range("L136").Value = 81624.5
$ws.Range("N136").Value = -91824.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 16666
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9826
$ws.Range("H58").Value = 2602.2
$ws.Range("I58").Value = 2196.05
$ws.Range("J58").Value = 3414.5
$ws.Range("K58").Value = 2196.05
$ws.Range("L58").Value = 3414.5
$ws.Range("M58").Value = -1993.05
$ws.Range("N58").Value = -3820.5
$ws.Range("H132").Value = 1877.1305
$ws.Range("I132").Value = 1894.2727
$ws.Range("K132").Value = 5682.8181
$ws.Range("M132").Value = -3152.8181
$ws.Range("H136").Value = 2602.2
$ws.Range("I136").Value = 2196.05
$ws.Range("J136").Value = 3414.5
$ws.Range("K136").Value = 6588.150000000001
$ws.Range("L136").Value = 10243.5
$ws.Range("M136").Value = -4038.150000000001
$ws.Range("N136").Value = -15343.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 298.0625
$ws.Range("I7").Value = 196.125
$ws.Range("K7").Value = 588.375
$ws.Range("M7").Value = -476.375
$ws.Range("H64").Value = 3997.8572
$ws.Range("I64").Value = 3997.8572
$ws.Range("K64").Value = 11993.5716
$ws.Range("M64").Value = -11723.5716
$ws.Range("H67").Value = 3997.8572
$ws.Range("I67").Value = 3997.8572
$ws.Range("K67").Value = 11993.5716
$ws.Range("M67").Value = -11057.5716
$ws.Range("H93").Value = 3649.6
$ws.Range("J93").Value = 4282.8335
$ws.Range("L93").Value = 12848.5005
$ws.Range("N93").Value = -16592.5005
$ws.Range("H107").Value = 20834264
$ws.Range("I107").Value = 33333810
$ws.Range("J107").Value = 1689.1666
$ws.Range("K107").Value = 100001430
$ws.Range("L107").Value = 5067.4998
$ws.Range("M107").Value = -99999510
$ws.Range("N107").Value = -8907.4998
$ws.Range("H108").Value = 398.7647
$ws.Range("I108").Value = 331.93332
$ws.Range("J108").Value = 900
$ws.Range("K108").Value = 995.7999599999999
$ws.Range("L108").Value = 2700
$ws.Range("M108").Value = 1884.20004
$ws.Range("N108").Value = -8460
$ws.Range("H119").Value = 804.75
$ws.Range("I119").Value = 804.75
$ws.Range("K119").Value = 2414.25
$ws.Range("M119").Value = 2423.75
$ws.Range("H122").Value = 5721
$ws.Range("J122").Value = 8335
$ws.Range("L122").Value = 75015
$ws.Range("N122").Value = -79915
$ws.Range("H140").Value = 3588.2942
$ws.Range("I140").Value = 3218.2727
$ws.Range("K140").Value = 9654.8181
$ws.Range("M140").Value = -4474.8181

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2392.5217
$ws.Range("I97").Value = 2264.3684
$ws.Range("K97").Value = 2264.3684
$ws.Range("M97").Value = -1768.3684
$ws.Range("H102").Value = 2439.1667
$ws.Range("I102").Value = 2160.9092
$ws.Range("K102").Value = 2160.9092
$ws.Range("M102").Value = -538.9092000000001
$ws.Range("H107").Value = 1164.4286
$ws.Range("I107").Value = 941.8333
$ws.Range("K107").Value = 941.8333
$ws.Range("M107").Value = 978.1667
$ws.Range("H122").Value = 16841.895
$ws.Range("I122").Value = 19350
$ws.Range("J122").Value = 7436.5
$ws.Range("K122").Value = 58050
$ws.Range("L122").Value = 22309.5
$ws.Range("M122").Value = -55600
$ws.Range("N122").Value = -27209.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 16142.128
$ws.Range("I132").Value = 18318.176
$ws.Range("J132").Value = 1345
$ws.Range("K132").Value = 54954.528
$ws.Range("L132").Value = 4035
$ws.Range("M132").Value = -52424.528
$ws.Range("N132").Value = -9095
$ws.Range("H135").Value = 48500
$ws.Range("J135").Value = 48500
$ws.Range("L135").Value = 48500
$ws.Range("N135").Value = -58640
$ws.Range("H136").Value = 44891
$ws.Range("J136").Value = 44891
$ws.Range("L136").Value = 134673
$ws.Range("N136").Value = -139773

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2670.9333
$ws.Range("I16").Value = 2524.423
$ws.Range("J16").Value = 3623.25
$ws.Range("K16").Value = 2524.423
$ws.Range("L16").Value = 3623.25
$ws.Range("M16").Value = -2354.423
$ws.Range("N16").Value = -3963.25
$ws.Range("H22").Value = 3127.5715
$ws.Range("I22").Value = 2662.25
$ws.Range("J22").Value = 3748
$ws.Range("K22").Value = 2662.25
$ws.Range("L22").Value = 3748
$ws.Range("M22").Value = -2367.25
$ws.Range("N22").Value = -4338
$ws.Range("H27").Value = 3127.5715
$ws.Range("I27").Value = 2662.25
$ws.Range("J27").Value = 3748
$ws.Range("K27").Value = 2662.25
$ws.Range("L27").Value = 3748
$ws.Range("M27").Value = -2555.25
$ws.Range("N27").Value = -3962
$ws.Range("H82").Value = 3140.4375
$ws.Range("I82").Value = 2508.2
$ws.Range("J82").Value = 4194.1665
$ws.Range("K82").Value = 2508.2
$ws.Range("L82").Value = 4194.1665
$ws.Range("M82").Value = -2147.2
$ws.Range("N82").Value = -4916.1665
$ws.Range("H85").Value = 3140.4375
$ws.Range("I85").Value = 2508.2
$ws.Range("J85").Value = 4194.1665
$ws.Range("K85").Value = 2508.2
$ws.Range("L85").Value = 4194.1665
$ws.Range("M85").Value = -1260.2
$ws.Range("N85").Value = -6690.1665
$ws.Range("H100").Value = 3960
$ws.Range("I100").Value = 3945
$ws.Range("K100").Value = 3945
$ws.Range("M100").Value = -3404
$ws.Range("H122").Value = 6985.5
$ws.Range("I122").Value = 6176.8
$ws.Range("K122").Value = 18530.4
$ws.Range("M122").Value = -16080.4
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 5126.091
$ws.Range("I132").Value = 4378.6
$ws.Range("K132").Value = 13135.8
$ws.Range("M132").Value = -10605.8

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 20000
$ws.Range("K43").Value = 20000
$ws.Range("M43").Value = -19851
$ws.Range("H46").Value = 71054.664
$ws.Range("J46").Value = 71054.664
$ws.Range("L46").Value = 71054.664
$ws.Range("N46").Value = -71516.664
$ws.Range("H62").Value = 18281.096
$ws.Range("I62").Value = 18209.162
$ws.Range("J62").Value = 18624.777
$ws.Range("K62").Value = 18209.162
$ws.Range("L62").Value = 18624.777
$ws.Range("M62").Value = -17585.162
$ws.Range("N62").Value = -19872.777
$ws.Range("H65").Value = 18281.096
$ws.Range("I65").Value = 18209.162
$ws.Range("J65").Value = 18624.777
$ws.Range("K65").Value = 91045.81
$ws.Range("L65").Value = 93123.88499999999
$ws.Range("M65").Value = -87925.81
$ws.Range("N65").Value = -99363.88499999999
$ws.Range("H101").Value = 23260.6
$ws.Range("J101").Value = 23260.6
$ws.Range("L101").Value = 23260.6
$ws.Range("N101").Value = -29750.6
$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 15000
$ws.Range("N104").Value = -21988
$ws.Range("H122").Value = 5237.645
$ws.Range("I122").Value = 1870.2858
$ws.Range("K122").Value = 5610.857400000001
$ws.Range("M122").Value = -3160.857400000001
$ws.Range("H125").Value = 500035740
$ws.Range("J125").Value = 500035740
$ws.Range("L125").Value = 500035740
$ws.Range("N125").Value = -500045580
$ws.Range("H126").Value = 3824.647
$ws.Range("I126").Value = 3201.2666
$ws.Range("K126").Value = 9603.799800000001
$ws.Range("M126").Value = -7133.799800000001
$ws.Range("H132").Value = 5179.095
$ws.Range("I132").Value = 6341.9287
$ws.Range("K132").Value = 19025.7861
$ws.Range("M132").Value = -16495.7861
$ws.Range("H134").Value = 71054.664
$ws.Range("J134").Value = 71054.664
$ws.Range("L134").Value = 213163.992
$ws.Range("N134").Value = -218233.992
$ws.Range("H136").Value = 2808.4
$ws.Range("I136").Value = 2239.8
$ws.Range("J136").Value = 3377
$ws.Range("K136").Value = 6719.400000000001
$ws.Range("L136").Value = 10131
$ws.Range("M136").Value = -4169.400000000001
$ws.Range("N136").Value = -15231
